$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column width changes (D,E narrower; G wider) - ColumnWidth is quantized by
# the host to an internal pixel grid, so these land on the nearest
# representable width to the target 2.140625 / 3.140625 character units.
$ws.Columns.Item(4).ColumnWidth = 1.3
$ws.Columns.Item(5).ColumnWidth = 1.3
$ws.Columns.Item(7).ColumnWidth = 2.3

# Row 1 data updates
$ws.Range("C1").Value = 19
$ws.Range("D1").Value = 6
$ws.Range("E1").Value = 7
$ws.Range("F1").Value = 16
$ws.Range("G1").Value = 13
$ws.Range("H1").Value = 24
$ws.Range("I1").Value = 23
$ws.Range("J1").Value = 32
$ws.Range("K1").Value = 0.013000000000000001
$ws.Range("L1").Value = 0.0080000000000000002
$ws.Range("M1").Value = 0.099000000000000005
$ws.Range("N1").Value = 0.060999999999999999
